$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 50.70817566666667
$ws.Range("H2").Value = 152.124527
$ws.Range("I2").Value = 0.5661129211027078
$ws.Range("J2").Value = 0.5661129211027077
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 2.237200333333333
$ws.Range("N2").Value = 6.711601
$ws.Range("O2").Value = 0.1121050933480713
$ws.Range("P2").Value = 0.1121050933480713
$ws.Range("Q2").Value = 113.4443475041919
$ws.Range("R2").Value = 1020.999127537727
$ws.Range("S2").Value = 0.0634641418657684
$ws.Range("T2").Value = 0.06346414186576839

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 50.70817566666667
$ws.Range("H3").Value = 152.124527
$ws.Range("I3").Value = 0.5661129211027078
$ws.Range("J3").Value = 0.5661129211027077
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 3.153682
$ws.Range("N3").Value = 9.461046
$ws.Range("O3").Value = 0.158029573718759
$ws.Range("P3").Value = 0.158029573718759
$ws.Range("Q3").Value = 159.9174608528047
$ws.Range("R3").Value = 1439.257147675242
$ws.Range("S3").Value = 0.08946258359854237
$ws.Range("T3").Value = 0.08946258359854237

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 50.70817566666667
$ws.Range("H4").Value = 152.124527
$ws.Range("I4").Value = 0.5661129211027078
$ws.Range("J4").Value = 0.5661129211027077
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 13.44189533333333
$ws.Range("N4").Value = 40.325686
$ws.Range("O4").Value = 0.6735672745377762
$ws.Range("P4").Value = 0.6735672745377762
$ws.Range("Q4").Value = 681.6139898556137
$ws.Range("R4").Value = 6134.525908700522
$ws.Range("S4").Value = 0.3813151373477701
$ws.Range("T4").Value = 0.38131513734777

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 50.70817566666667
$ws.Range("H5").Value = 152.124527
$ws.Range("I5").Value = 0.5661129211027078
$ws.Range("J5").Value = 0.5661129211027077
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 1.123499666666667
$ws.Range("N5").Value = 3.370499
$ws.Range("O5").Value = 0.05629805839539345
$ws.Range("P5").Value = 0.05629805839539345
$ws.Range("Q5").Value = 56.97061845877479
$ws.Range("R5").Value = 512.735566128973
$ws.Range("S5").Value = 0.03187105829062701
$ws.Range("T5").Value = 0.031871058290627

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 17.08683666666667
$ws.Range("H6").Value = 51.26051
$ws.Range("I6").Value = 0.1907597520636141
$ws.Range("J6").Value = 0.1907597520636141
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 2.237200333333333
$ws.Range("N6").Value = 6.711601
$ws.Range("O6").Value = 0.1121050933480713
$ws.Range("P6").Value = 0.1121050933480713
$ws.Range("Q6").Value = 38.22667668627889
$ws.Range("R6").Value = 344.04009017651
$ws.Range("S6").Value = 0.0213851398121464
$ws.Range("T6").Value = 0.0213851398121464

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 17.08683666666667
$ws.Range("H7").Value = 51.26051
$ws.Range("I7").Value = 0.1907597520636141
$ws.Range("J7").Value = 0.1907597520636141
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 3.153682
$ws.Range("N7").Value = 9.461046
$ws.Range("O7").Value = 0.158029573718759
$ws.Range("P7").Value = 0.158029573718759
$ws.Range("Q7").Value = 53.88644923260667
$ws.Range("R7").Value = 484.97804309346
$ws.Range("S7").Value = 0.03014568230130909
$ws.Range("T7").Value = 0.03014568230130909

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 17.08683666666667
$ws.Range("H8").Value = 51.26051
$ws.Range("I8").Value = 0.1907597520636141
$ws.Range("J8").Value = 0.1907597520636141
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 13.44189533333333
$ws.Range("N8").Value = 40.325686
$ws.Range("O8").Value = 0.6735672745377762
$ws.Range("P8").Value = 0.6735672745377762
$ws.Range("Q8").Value = 229.6794700510956
$ws.Range("R8").Value = 2067.11523045986
$ws.Range("S8").Value = 0.1284895262889905
$ws.Range("T8").Value = 0.1284895262889905

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 17.08683666666667
$ws.Range("H9").Value = 51.26051
$ws.Range("I9").Value = 0.1907597520636141
$ws.Range("J9").Value = 0.1907597520636141
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 1.123499666666667
$ws.Range("N9").Value = 3.370499
$ws.Range("O9").Value = 0.05629805839539345
$ws.Range("P9").Value = 0.05629805839539345
$ws.Range("Q9").Value = 19.19705529938778
$ws.Range("R9").Value = 172.77349769449
$ws.Range("S9").Value = 0.01073940366116812
$ws.Range("T9").Value = 0.01073940366116812

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 18.33915266666667
$ws.Range("H10").Value = 55.017458
$ws.Range("I10").Value = 0.2047407770084672
$ws.Range("J10").Value = 0.2047407770084672
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 2.237200333333333
$ws.Range("N10").Value = 6.711601
$ws.Range("O10").Value = 0.1121050933480713
$ws.Range("P10").Value = 0.1121050933480713
$ws.Range("Q10").Value = 41.02835845891756
$ws.Range("R10").Value = 369.2552261302581
$ws.Range("S10").Value = 0.02295248391869087
$ws.Range("T10").Value = 0.02295248391869087

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 18.33915266666667
$ws.Range("H11").Value = 55.017458
$ws.Range("I11").Value = 0.2047407770084672
$ws.Range("J11").Value = 0.2047407770084672
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 3.153682
$ws.Range("N11").Value = 9.461046
$ws.Range("O11").Value = 0.158029573718759
$ws.Range("P11").Value = 0.158029573718759
$ws.Range("Q11").Value = 57.83585566011867
$ws.Range("R11").Value = 520.522700941068
$ws.Range("S11").Value = 0.03235509771349556
$ws.Range("T11").Value = 0.03235509771349557

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 18.33915266666667
$ws.Range("H12").Value = 55.017458
$ws.Range("I12").Value = 0.2047407770084672
$ws.Range("J12").Value = 0.2047407770084672
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 13.44189533333333
$ws.Range("N12").Value = 40.325686
$ws.Range("O12").Value = 0.6735672745377762
$ws.Range("P12").Value = 0.6735672745377762
$ws.Range("Q12").Value = 246.5129706473543
$ws.Range("R12").Value = 2218.616735826188
$ws.Range("S12").Value = 0.1379066871563398
$ws.Range("T12").Value = 0.1379066871563398

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 18.33915266666667
$ws.Range("H13").Value = 55.017458
$ws.Range("I13").Value = 0.2047407770084672
$ws.Range("J13").Value = 0.2047407770084672
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 1.123499666666667
$ws.Range("N13").Value = 3.370499
$ws.Range("O13").Value = 0.05629805839539345
$ws.Range("P13").Value = 0.05629805839539345
$ws.Range("Q13").Value = 20.60403190794911
$ws.Range("R13").Value = 185.436287171542
$ws.Range("S13").Value = 0.01152650821994091
$ws.Range("T13").Value = 0.01152650821994091

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 3.438381
$ws.Range("H14").Value = 10.315143
$ws.Range("I14").Value = 0.03838654982521095
$ws.Range("J14").Value = 0.03838654982521095
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 2.237200333333333
$ws.Range("N14").Value = 6.711601
$ws.Range("O14").Value = 0.1121050933480713
$ws.Range("P14").Value = 0.1121050933480713
$ws.Range("Q14").Value = 7.692347119327001
$ws.Range("R14").Value = 69.23112407394301
$ws.Range("S14").Value = 0.004303327751465665
$ws.Range("T14").Value = 0.004303327751465665

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 3.438381
$ws.Range("H15").Value = 10.315143
$ws.Range("I15").Value = 0.03838654982521095
$ws.Range("J15").Value = 0.03838654982521095
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 3.153682
$ws.Range("N15").Value = 9.461046
$ws.Range("O15").Value = 0.158029573718759
$ws.Range("P15").Value = 0.158029573718759
$ws.Range("Q15").Value = 10.843560268842
$ws.Range("R15").Value = 97.59204241957801
$ws.Range("S15").Value = 0.00606621010541199
$ws.Range("T15").Value = 0.006066210105411991

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 3.438381
$ws.Range("H16").Value = 10.315143
$ws.Range("I16").Value = 0.03838654982521095
$ws.Range("J16").Value = 0.03838654982521095
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 13.44189533333333
$ws.Range("N16").Value = 40.325686
$ws.Range("O16").Value = 0.6735672745377762
$ws.Range("P16").Value = 0.6735672745377762
$ws.Range("Q16").Value = 46.218357518122
$ws.Range("R16").Value = 415.9652176630981
$ws.Range("S16").Value = 0.02585592374467589
$ws.Range("T16").Value = 0.02585592374467589

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 3.438381
$ws.Range("H17").Value = 10.315143
$ws.Range("I17").Value = 0.03838654982521095
$ws.Range("J17").Value = 0.03838654982521095
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 1.123499666666667
$ws.Range("N17").Value = 3.370499
$ws.Range("O17").Value = 0.05629805839539345
$ws.Range("P17").Value = 0.05629805839539345
$ws.Range("Q17").Value = 3.863019907373001
$ws.Range("R17").Value = 34.767179166357
$ws.Range("S17").Value = 0.002161088223657406
$ws.Range("T17").Value = 0.002161088223657406
